$wb = $excel.ActiveWorkbook

$negs = $wb.Worksheets.Item("Negs")
$negs.Range("E1").Value = "Reason"
$negs.Range("E1").Font.Bold = $true
$negs.Range("E2").Value = "Do not want to work"

$leave = $wb.Worksheets.Item("Leave")
$leave.Range("D1").Value = "Reason"
$leave.Range("D1").Font.Bold = $true
$leave.Range("D2").Value = "Holiday"

$leave.Activate()
